$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the affected rows
$ws.Range("F2").Value = -9
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = -7
$ws.Range("F8").Value = -1
